$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.382.11"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "3.340.95"
$ws.Range("E3").Value = "  -3.50%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'551.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.56%  "

$ws.Range("D6").Value = "'174.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "

$ws.Range("D7").Value = "'0.618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.18%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.338.00"
$ws.Range("E8").Value = "  -3.27%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("D12").Value = "'54.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("E13").Value = "  -1.57%  "

$ws.Range("D14").Value = "'9.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").Value = "3.873.44"
$ws.Range("E15").Value = "  -3.63%  "

$ws.Range("D16").Value = "'18.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("E17").Value = "  -2.95%  "

$ws.Range("D18").Value = "3.337.98"
$ws.Range("E18").Value = "  -3.79%  "

$ws.Range("D19").Value = "64.291.28"
$ws.Range("E19").Value = "  -1.84%  "

$ws.Range("D20").Value = "'11.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.56%  "

$ws.Range("E21").Value = "  -2.62%  "

$ws.Range("D22").Value = "'438.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.91%  "

$ws.Range("D23").Value = "'5.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.80%  "

$ws.Range("D24").Value = "'4.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.60%  "

$ws.Range("D25").Value = "'84.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "'13.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").Value = "'10.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.22%  "

$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("E29").Value = "  -4.69%  "

$ws.Range("E30").Value = "  -1.38%  "

$ws.Range("E31").Value = "  -1.26%  "

$ws.Range("D32").Value = "'11.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.10%  "

$ws.Range("D33").Value = "'576.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.82%  "

$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("D35").Value = "'58.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.84%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -8.16%  "

$ws.Range("D38").Value = "'3.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("D39").Value = "'35.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.58%  "

$ws.Range("D40").Value = "0.0₃0747"
$ws.Range("E40").Value = "  -5.24%  "

$ws.Range("E41").Value = "  -4.11%  "

$ws.Range("D42").Value = "3.097.68"
$ws.Range("E42").Value = "  -3.75%  "

$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("E44").Value = "  -6.25%  "

$ws.Range("D45").Value = "'3.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.98%  "

$ws.Range("E46").Value = "  -1.95%  "

$ws.Range("E47").Value = "  -3.33%  "

$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("E49").Value = "  -1.57%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'136.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.41%  "
